# Update countries & provincias Spain
#
# Applies the daily COVID data refresh to the "Pais" sheet:
#   - bumps the "last updated" timestamp in A1
#   - reorders two pairs/blocks in the shared country list:
#       * Eritrea now sorts immediately before Mongolia (they swap places)
#       * Lesoto now sorts immediately after Brunei / before Trinidad yTobago
#         (Trinidad yTobago, Monaco, Bahamas, Aruba, Barbados each shift down
#         one row; Seychelles keeps its row)
#   - refreshes the per-country case numbers (columns B:H) for every row
#     whose figures changed in this update

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados a ..." banner -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Julio de 2020 a las 01:58"

# --- Reordering within the country list (column A) -------------------------
# Eritrea / Mongolia swap (rows 169-170)
$ws.Range("A169").Value = "Eritrea"
$ws.Range("A170").Value = "Mongolia"

# Lesoto inserted right after Brunei; everything through Barbados shifts down
# one row (rows 179-184). Seychelles (row 185) keeps its position.
$ws.Range("A179").Value = "Lesoto"
$ws.Range("A180").Value = "Trinidad yTobago"
$ws.Range("A181").Value = "Monaco"
$ws.Range("A182").Value = "Bahamas"
$ws.Range("A183").Value = "Aruba"
$ws.Range("A184").Value = "Barbados"

# --- Refreshed statistics (columns B:H) -------------------------------------
# Each entry: row, Casos totales, Nuevos casos, Casos activos, Recuperados,
#             Casos criticos, Muertes hoy, Muertes
$updates = @(
    @(4,   3218903, 59971, 1426170, 1656939, 0, 932, 135794),
    @(5,   1759103, 42907, 1152467,  537382, 0, 1199,  69254),
    @(25,    90693,  3663,   38313,   50660, 0,   26,   1720),
    @(75,     8965,    15,    8138,     575, 0,    1,    252),
    @(124,    1598,    14,    1123,     412, 0,    0,     63),
    @(142,     977,     3,     878,      70, 0,    0,     29),
    @(149,     726,     2,     283,     429, 0,    1,     14),
    @(152,     699,     1,     656,       1, 0,    0,     42),
    @(166,     314,     1,     272,      35, 0,    0,      7),
    @(167,     286,     2,     125,     145, 0,    0,     16),
    @(169,     232,    17,     107,     125, 0,    0,      0),
    @(170,     227,     0,     197,      30, 0,    0,      0),
    @(179,     134,    43,      11,     122, 0,    1,      1),
    @(180,     133,     0,     120,       5, 0,    0,      8),
    @(181,     108,     0,      96,       8, 0,    0,      4),
    @(182,     106,     2,      89,       6, 0,    0,     11),
    @(183,     105,     0,      98,       4, 0,    0,      3),
    @(184,      98,     0,      90,       1, 0,    0,      7),
    @(185,      94,     3,      11,      83, 0,    0,      0)
)

foreach ($row in $updates) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]   # B - Casos totales
    $ws.Cells.Item($r, 3).Value = $row[2]   # C - Nuevos casos
    $ws.Cells.Item($r, 4).Value = $row[3]   # D - Casos activos
    $ws.Cells.Item($r, 5).Value = $row[4]   # E - Recuperados
    $ws.Cells.Item($r, 6).Value = $row[5]   # F - Casos criticos
    $ws.Cells.Item($r, 7).Value = $row[6]   # G - Muertes hoy
    $ws.Cells.Item($r, 8).Value = $row[7]   # H - Muertes
}
